$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '69.160.08'
$ws.Range('E2').Value = '  +2.01%  '
$ws.Range('D3').Value = '3.833.40'
$ws.Range('E3').Value = '  +0.59%  '
$ws.Range('D4').Value = "'0.999"
$ws.Range('D4').ClearFormats()
$ws.Range('E4').Value = '  -0.72%  '
$ws.Range('D5').Value = "'630.04"
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +4.69%  '
$ws.Range('D6').Value = "'165.51"
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -0.43%  '
$ws.Range('D7').Value = '3.830.97'
$ws.Range('E7').Value = '  +0.62%  '
$ws.Range('E8').Value = '  -0.05%  '
$ws.Range('D9').Value = "'0.520"
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.27%  '
$ws.Range('E10').Value = '  +1.26%  '
$ws.Range('D11').Value = "'0.455"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +0.34%  '
$ws.Range('D12').Value = "'6.72"
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +5.89%  '
$ws.Range('E13').Value = '  +0.21%  '
$ws.Range('D14').Value = "'35.92"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -0.32%  '
$ws.Range('D15').Value = '4.476.24'
$ws.Range('E15').Value = '  +0.61%  '
$ws.Range('D16').Value = '3.827.19'
$ws.Range('E16').Value = '  +0.67%  '
$ws.Range('D17').Value = '69.162.38'
$ws.Range('E17').Value = '  +1.95%  '
$ws.Range('D18').Value = "'18.24"
$ws.Range('D18').ClearFormats()
$ws.Range('E18').Value = '  -1.29%  '
$ws.Range('D19').Value = "'7.16"
$ws.Range('D19').ClearFormats()
$ws.Range('E19').Value = '  +0.98%  '
$ws.Range('E20').Value = '  -0.04%  '
$ws.Range('D21').Value = "'467.94"
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  +0.86%  '
$ws.Range('D22').Value = "'9.77"
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -0.53%  '
$ws.Range('D23').Value = "'0.709"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +0.76%  '
$ws.Range('E24').Value = '  +3.57%  '
$ws.Range('D25').Value = "'84.08"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  +1.02%  '
$ws.Range('D26').Value = "'12.10"
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.10%  '
$ws.Range('E27').Value = '  +1.98%  '
$ws.Range('D28').Value = "'10.08"
$ws.Range('D28').ClearFormats()
$ws.Range('E28').Value = '  +0.67%  '
$ws.Range('E29').Value = '  +0.04%  '
$ws.Range('D30').Value = '3.981.85'
$ws.Range('E30').Value = '  +0.55%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').Value = "'2.67"
$ws.Range('D31').ClearFormats()
$ws.Range('E31').Value = '  -3.81%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = "'2.25"
$ws.Range('D32').ClearFormats()
$ws.Range('E32').Value = '  +1.80%  '
$ws.Range('D33').Value = "'7.33"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  -1.23%  '
$ws.Range('D34').Value = "'29.25"
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -0.59%  '
$ws.Range('D35').Value = "'9.13"
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  +0.58%  '
$ws.Range('E37').Value = '  +2.09%  '
$ws.Range('E38').Value = '  +7.08%  '
$ws.Range('D39').Value = "'5.93"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  +1.93%  '
$ws.Range('E40').Value = '  +1.39%  '
$ws.Range('D41').Value = "'0.980"
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  -1.61%  '
$ws.Range('D42').Value = "'0.999"
$ws.Range('D42').ClearFormats()
$ws.Range('E42').Value = '  -0.09%  '
$ws.Range('D44').Value = "'157.20"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  +3.94%  '
$ws.Range('E45').Value = '  +0.62%  '
$ws.Range('D46').Value = "'1.42"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  +3.33%  '
$ws.Range('D47').Value = "'43.02"
$ws.Range('D47').ClearFormats()
$ws.Range('E47').Value = '  -4.80%  '
$ws.Range('D48').Value = "'46.92"
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  -1.78%  '
$ws.Range('E49').Value = '  +2.73%  '
$ws.Range('D50').Value = "'8.43"
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  +0.84%  '
$ws.Range('D51').Value = "'382.60"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.24%  '
